$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Paragraphs.Item(1).Range.Text = "2024-08-05 Monday"

# Update each table cell's equation text in row-major order
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "86-24=62"  # was: 50+0=50
$t.Cell(1, 2).Range.Text = "34+32=66"  # was: 71-17=54
$t.Cell(1, 3).Range.Text = "6+72=78"  # was: 66+6=72
$t.Cell(1, 4).Range.Text = "22+7=29"  # was: 89-71=18
$t.Cell(1, 5).Range.Text = "5+20=25"  # was: 15+70=85
$t.Cell(2, 1).Range.Text = "86-51=35"  # was: 4+58=62
$t.Cell(2, 2).Range.Text = "79-74=5"  # was: 30-27=3
$t.Cell(2, 3).Range.Text = "94-13=81"  # was: 32-27=5
$t.Cell(2, 4).Range.Text = "1+75=76"  # was: 64+13=77
$t.Cell(2, 5).Range.Text = "54-53=1"  # was: 68+23=91
$t.Cell(3, 1).Range.Text = "86-17=69"  # was: 29+69=98
$t.Cell(3, 2).Range.Text = "59+29=88"  # was: 6+77=83
$t.Cell(3, 3).Range.Text = "27+56=83"  # was: 32-8=24
$t.Cell(3, 4).Range.Text = "94-68=26"  # was: 58+3=61
$t.Cell(3, 5).Range.Text = "42+43=85"  # was: 69-54=15
$t.Cell(4, 1).Range.Text = "79+9=88"  # was: 88-31=57
$t.Cell(4, 2).Range.Text = "0+88=88"  # was: 84-75=9
$t.Cell(4, 3).Range.Text = "38+20=58"  # was: 82-54=28
$t.Cell(4, 4).Range.Text = "72-56=16"  # was: 38+19=57
$t.Cell(4, 5).Range.Text = "90-41=49"  # was: 50+18=68
$t.Cell(5, 1).Range.Text = "39+9=48"  # was: 0+36=36
$t.Cell(5, 2).Range.Text = "67-50=17"  # was: 11+76=87
$t.Cell(5, 3).Range.Text = "0+99=99"  # was: 71-48=23
$t.Cell(5, 4).Range.Text = "19+34=53"  # was: 40+47=87
$t.Cell(5, 5).Range.Text = "59-44=15"  # was: 51-1=50
$t.Cell(6, 1).Range.Text = "33-14=19"  # was: 83-52=31
$t.Cell(6, 2).Range.Text = "97-53=44"  # was: 19+10=29
$t.Cell(6, 3).Range.Text = "16+12=28"  # was: 73-0=73
$t.Cell(6, 4).Range.Text = "72-67=5"  # was: 78+4=82
$t.Cell(6, 5).Range.Text = "18+78=96"  # was: 29+14=43
$t.Cell(7, 1).Range.Text = "35+12=47"  # was: 57-32=25
$t.Cell(7, 2).Range.Text = "74-74=0"  # was: 15+68=83
$t.Cell(7, 3).Range.Text = "49-40=9"  # was: 43+40=83
$t.Cell(7, 4).Range.Text = "85-61=24"  # was: 64+30=94
$t.Cell(7, 5).Range.Text = "45+42=87"  # was: 59+38=97
$t.Cell(8, 1).Range.Text = "44-25=19"  # was: 28-15=13
$t.Cell(8, 2).Range.Text = "94-76=18"  # was: 92-91=1
$t.Cell(8, 3).Range.Text = "75+10=85"  # was: 13+67=80
$t.Cell(8, 4).Range.Text = "34+50=84"  # was: 97-15=82
$t.Cell(8, 5).Range.Text = "23+1=24"  # was: 76+5=81
$t.Cell(9, 1).Range.Text = "56-37=19"  # was: 26-2=24
$t.Cell(9, 2).Range.Text = "15+62=77"  # was: 35+60=95
$t.Cell(9, 3).Range.Text = "59-7=52"  # was: 46+22=68
$t.Cell(9, 4).Range.Text = "68+30=98"  # was: 54-24=30
$t.Cell(9, 5).Range.Text = "41+41=82"  # was: 9+71=80
$t.Cell(10, 1).Range.Text = "45+40=85"  # was: 18+8=26
$t.Cell(10, 2).Range.Text = "21+7=28"  # was: 77-55=22
$t.Cell(10, 3).Range.Text = "77-53=24"  # was: 72-4=68
$t.Cell(10, 4).Range.Text = "62-6=56"  # was: 2+41=43
$t.Cell(10, 5).Range.Text = "30+39=69"  # was: 64+35=99
$t.Cell(11, 1).Range.Text = "65-60=5"  # was: 95-48=47
$t.Cell(11, 2).Range.Text = "34+6=40"  # was: 32+18=50
$t.Cell(11, 3).Range.Text = "48-29=19"  # was: 30+65=95
$t.Cell(11, 4).Range.Text = "15-10=5"  # was: 58+32=90
$t.Cell(11, 5).Range.Text = "41+26=67"  # was: 80+10=90
$t.Cell(12, 1).Range.Text = "32-5=27"  # was: 94-71=23
$t.Cell(12, 2).Range.Text = "26+17=43"  # was: 86-65=21
$t.Cell(12, 3).Range.Text = "15+18=33"  # was: 46-14=32
$t.Cell(12, 4).Range.Text = "71-22=49"  # was: 40-24=16
$t.Cell(12, 5).Range.Text = "25-6=19"  # was: 53-2=51
$t.Cell(13, 1).Range.Text = "2+18=20"  # was: 33-11=22
$t.Cell(13, 2).Range.Text = "9+41=50"  # was: 80-77=3
$t.Cell(13, 3).Range.Text = "84-72=12"  # was: 7+28=35
$t.Cell(13, 4).Range.Text = "62-43=19"  # was: 94-7=87
$t.Cell(13, 5).Range.Text = "38-22=16"  # was: 90-1=89
$t.Cell(14, 1).Range.Text = "8+75=83"  # was: 66+0=66
$t.Cell(14, 2).Range.Text = "51-2=49"  # was: 82-4=78
$t.Cell(14, 3).Range.Text = "0+78=78"  # was: 25+19=44
$t.Cell(14, 4).Range.Text = "52+24=76"  # was: 22+33=55
$t.Cell(14, 5).Range.Text = "60-41=19"  # was: 5+79=84
$t.Cell(15, 1).Range.Text = "47+10=57"  # was: 32-19=13
$t.Cell(15, 2).Range.Text = "62+18=80"  # was: 15+1=16
$t.Cell(15, 3).Range.Text = "17+19=36"  # was: 94-24=70
$t.Cell(15, 4).Range.Text = "67-61=6"  # was: 7+37=44
$t.Cell(15, 5).Range.Text = "28-15=13"  # was: 28+67=95
$t.Cell(16, 1).Range.Text = "56-12=44"  # was: 59+38=97
$t.Cell(16, 2).Range.Text = "82-71=11"  # was: 25+5=30
$t.Cell(16, 3).Range.Text = "48-35=13"  # was: 48-4=44
$t.Cell(16, 4).Range.Text = "51+22=73"  # was: 49-45=4
$t.Cell(16, 5).Range.Text = "26-18=8"  # was: 24+29=53
$t.Cell(17, 1).Range.Text = "80-37=43"  # was: 94-32=62
$t.Cell(17, 2).Range.Text = "93-76=17"  # was: 92-22=70
$t.Cell(17, 3).Range.Text = "69-39=30"  # was: 91-60=31
$t.Cell(17, 4).Range.Text = "16+49=65"  # was: 3+79=82
$t.Cell(17, 5).Range.Text = "45+40=85"  # was: 54-0=54
$t.Cell(18, 1).Range.Text = "30-21=9"  # was: 24+60=84
$t.Cell(18, 2).Range.Text = "36+62=98"  # was: 92-51=41
$t.Cell(18, 3).Range.Text = "29+65=94"  # was: 56+21=77
$t.Cell(18, 4).Range.Text = "47-15=32"  # was: 91-79=12
$t.Cell(18, 5).Range.Text = "43-30=13"  # was: 97-69=28
$t.Cell(19, 1).Range.Text = "82-33=49"  # was: 78-69=9
$t.Cell(19, 2).Range.Text = "33-27=6"  # was: 42+34=76
$t.Cell(19, 3).Range.Text = "78-28=50"  # was: 6+90=96
$t.Cell(19, 4).Range.Text = "71-33=38"  # was: 69-24=45
$t.Cell(19, 5).Range.Text = "47+12=59"  # was: 15+64=79
$t.Cell(20, 1).Range.Text = "56+13=69"  # was: 13-11=2
$t.Cell(20, 2).Range.Text = "73-50=23"  # was: 29+39=68
$t.Cell(20, 3).Range.Text = "75-1=74"  # was: 59-41=18
$t.Cell(20, 4).Range.Text = "50-48=2"  # was: 27+72=99
$t.Cell(20, 5).Range.Text = "52-30=22"  # was: 95-75=20
